$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for column G ("K") rows 2-39, replacing old Strike# values
# with recalculated K values (regen save_data to use K instead of Strike#)
$gValues = @{
    2  = 0
    3  = 8
    4  = 3
    5  = 5
    6  = 2
    7  = 5
    8  = 5
    9  = 5
    10 = 4
    11 = 5
    12 = 4
    13 = 7
    14 = 3
    15 = 1
    16 = 8
    17 = 10
    18 = 5
    19 = 5
    20 = 7
    21 = 6
    22 = 5
    23 = 6
    24 = 5
    25 = 8
    26 = 4
    27 = 4
    28 = 2
    29 = 6
    30 = 5
    31 = 6
    32 = 1
    33 = 3
    34 = 3
    35 = 7
    36 = 2
    37 = 3
    38 = 3
    39 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
